$d = $word.ActiveDocument

# --- Change 1: append "  (This is a change – Version for branch alternate)"
#     to the end of the first paragraph; the "(This is a change ... alternate)"
#     portion is coloured dark red (C00000).
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.InsertAfter("  ")

# Paragraph.Range.End points one past the trailing paragraph mark, so
# back up one position to get the real "end of visible text" offset
# before appending each new fragment.
$dash = [char]0x2013

$r2 = $d.Paragraphs(1).Range
$s1Start = $r2.End - 1
$r2.InsertAfter("(This is a change ${dash} Ve")
$s1End = $d.Paragraphs(1).Range.End - 1
$d.Range($s1Start, $s1End).Font.Color = 192

$r3 = $d.Paragraphs(1).Range
$s2Start = $r3.End - 1
$r3.InsertAfter("rsion for branch alternate")
$s2End = $d.Paragraphs(1).Range.End - 1
$d.Range($s2Start, $s2End).Font.Color = 192

$r4 = $d.Paragraphs(1).Range
$s3Start = $r4.End - 1
$r4.InsertAfter(")")
$s3End = $d.Paragraphs(1).Range.End - 1
$d.Range($s3Start, $s3End).Font.Color = 192

# --- Change 2: add a new, blank paragraph at the very end of the body
#     (before the section break).
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$newLast = $d.Paragraphs($d.Paragraphs.Count)
$newLast.Style = "Normal"
